# Updates cryptos list values (price/volume) per the latest scrape.
# Column D holds price text like "29.121.13" or "1.002" that must stay
# TEXT (not be reinterpreted as a number/date) -- exactly as it was originally
# authored as an inline string. A leading apostrophe forces Excel to keep the
# cell as text, matching the source workbook's inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.121.13'
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").Value = '1.835.70'
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").Value = '''244.05'
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("D6").Value = '''0.6284'
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").Value = '''1.002'
$ws.Range("E7").Value = '  +0.26%  '

$ws.Range("D8").Value = '''0.07538'
$ws.Range("E8").Value = '  -0.62%  '

$ws.Range("D9").Value = '''0.2927'
$ws.Range("E9").Value = '  +0.09%  '

$ws.Range("D10").Value = '''23.21'
$ws.Range("E10").Value = '  +2.72%  '

$ws.Range("D11").Value = '''0.07722'
$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("D12").Value = '1.825.99'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("D13").Value = '''5.011'
$ws.Range("E13").Value = '  +1.17%  '

$ws.Range("D14").Value = '''0.6693'
$ws.Range("E14").Value = '  +0.46%  '

$ws.Range("D15").Value = '''82.70'
$ws.Range("E15").Value = '  -0.03%  '

$ws.Range("D16").Value = '''0.000009424'
$ws.Range("E16").Value = '  -7.40%  '

$ws.Range("D17").Value = '''5.995'
$ws.Range("E17").Value = '  -0.75%  '

$ws.Range("D18").Value = '29.125.37'
$ws.Range("E18").Value = '  +0.32%  '

$ws.Range("D19").Value = '2.074.25'
$ws.Range("E19").Value = '  -0.55%  '

$ws.Range("D20").Value = '''12.60'
$ws.Range("E20").Value = '  +2.02%  '

$ws.Range("D21").Value = '''224.05'
$ws.Range("E21").Value = '  -1.17%  '

$ws.Range("D22").Value = '''1.005'
$ws.Range("E22").Value = '  +0.66%  '

$ws.Range("D23").Value = '''7.112'
$ws.Range("E23").Value = '  -0.90%  '

$ws.Range("D24").Value = '''1.003'
$ws.Range("E24").Value = '  +0.32%  '

$ws.Range("D25").Value = '''160.34'
$ws.Range("E25").Value = '  +1.25%  '

$ws.Range("D26").Value = '''0.1398'

$ws.Range("D27").Value = '''8.510'
$ws.Range("E27").Value = '  +0.17%  '

$ws.Range("D28").Value = '''17.96'
$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("D29").Value = '''1.501'
$ws.Range("E29").Value = '  +0.62%  '

$ws.Range("D30").Value = '''0.05657'
$ws.Range("E30").Value = '  +8.42%  '

$ws.Range("D31").Value = '''4.159'
$ws.Range("E31").Value = '  +1.35%  '

$ws.Range("D32").Value = '''4.068'
$ws.Range("E32").Value = '  +1.21%  '

$ws.Range("D33").Value = '''1.203'
$ws.Range("E33").Value = '  +0.95%  '

$ws.Range("D34").Value = '''0.7477'
$ws.Range("E34").Value = '  +1.47%  '

$ws.Range("D35").Value = '''1.849'
$ws.Range("E35").Value = '  +0.14%  '

$ws.Range("D36").Value = '''1.140'
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("D37").Value = '''2.673'
$ws.Range("E37").Value = '  -1.09%  '

$ws.Range("D38").Value = '''2.764'
$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("D39").Value = '1.223.13'
$ws.Range("E39").Value = '  -1.22%  '

$ws.Range("D40").Value = '''0.01784'
$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("D41").Value = '''6.517'
$ws.Range("E41").Value = '  +2.85%  '

$ws.Range("D42").Value = '''0.8949'
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").Value = '''1.003'
$ws.Range("E43").Value = '  +0.30%  '

$ws.Range("D44").Value = '''102.00'
$ws.Range("E44").Value = '  +0.47%  '

$ws.Range("D45").Value = '1.976.10'
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '''0.00000000125'
$ws.Range("E46").Value = '  +0.12%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''65.70'
$ws.Range("E47").Value = '  +2.21%  '

$ws.Range("B48").Value = 'XinFinNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D48").Value = '''0.07690'
$ws.Range("E48").Value = '  +11.86%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '''0.5094'
$ws.Range("E49").Value = '  -0.28%  '

$ws.Range("D50").Value = '''0.4082'
$ws.Range("E50").Value = '  +1.05%  '

$ws.Range("D51").Value = '''9.065'
$ws.Range("E51").Value = '  +2.42%  '
